$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3457.5386
$ws.Range("I62").Value = 3162.3333
$ws.Range("K62").Value = 3162.3333
$ws.Range("M62").Value = -2538.3333
$ws.Range("H65").Value = 3457.5386
$ws.Range("I65").Value = 3162.3333
$ws.Range("K65").Value = 15811.6665
$ws.Range("M65").Value = -12691.6665
$ws.Range("H69").Value = 4500
$ws.Range("I69").Value = 4500
$ws.Range("K69").Value = 13500
$ws.Range("M69").Value = -12626
$ws.Range("H72").Value = 4500
$ws.Range("I72").Value = 4500
$ws.Range("K72").Value = 40500
$ws.Range("M72").Value = -36132
$ws.Range("H74").Value = 2512.611
$ws.Range("I74").Value = 1895.4375
$ws.Range("J74").Value = 7450
$ws.Range("K74").Value = 1895.4375
$ws.Range("L74").Value = 7450
$ws.Range("M74").Value = -959.4375
$ws.Range("N74").Value = -9322
$ws.Range("H77").Value = 2512.611
$ws.Range("I77").Value = 1895.4375
$ws.Range("J77").Value = 7450
$ws.Range("K77").Value = 9477.1875
$ws.Range("L77").Value = 37250
$ws.Range("M77").Value = -4797.1875
$ws.Range("N77").Value = -46610
$ws.Range("H80").Value = 537.93335
$ws.Range("I80").Value = 471.45
$ws.Range("J80").Value = 670.9
$ws.Range("K80").Value = 1414.35
$ws.Range("L80").Value = 2012.7
$ws.Range("M80").Value = -416.3499999999999
$ws.Range("N80").Value = -4008.7
$ws.Range("H83").Value = 537.93335
$ws.Range("I83").Value = 471.45
$ws.Range("J83").Value = 670.9
$ws.Range("K83").Value = 4243.05
$ws.Range("L83").Value = 6038.099999999999
$ws.Range("M83").Value = 748.9499999999998
$ws.Range("N83").Value = -16022.1
$ws.Range("H106").Value = 2639.2
$ws.Range("I106").Value = 1999
$ws.Range("J106").Value = 3599.5
$ws.Range("K106").Value = 1999
$ws.Range("L106").Value = 3599.5
$ws.Range("M106").Value = -1368
$ws.Range("N106").Value = -4861.5
$ws.Range("H112").Value = 7191.514
$ws.Range("J112").Value = 9248.579
$ws.Range("L112").Value = 27745.737
$ws.Range("N112").Value = -29961.737
$ws.Range("H116").Value = 4237821.5
$ws.Range("I116").Value = 7412688.5
$ws.Range("K116").Value = 7412688.5
$ws.Range("M116").Value = -7409246.5
$ws.Range("H137").Value = 107007.95
$ws.Range("I137").Value = 2072
$ws.Range("J137").Value = 246922.56
$ws.Range("K137").Value = 6216
$ws.Range("L137").Value = 740767.6799999999
$ws.Range("M137").Value = -3666
$ws.Range("N137").Value = -745867.6799999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2504294.2
$ws.Range("I32").Value = 2859733.2
$ws.Range("K32").Value = 2859733.2
$ws.Range("M32").Value = -2859446.2
$ws.Range("H61").Value = 2433931.5
$ws.Range("I61").Value = 3257.4443
$ws.Range("K61").Value = 3257.4443
$ws.Range("M61").Value = -3045.4443
$ws.Range("H63").Value = 3225
$ws.Range("I63").Value = 3225
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3225
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -2539
$ws.Range("H66").Value = 3225
$ws.Range("I66").Value = 3225
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 16125
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -12693
$ws.Range("H97").Value = 287.23077
$ws.Range("I97").Value = 287.23077
$ws.Range("K97").Value = 287.23077
$ws.Range("M97").Value = 208.76923
$ws.Range("H102").Value = 4836.3335
$ws.Range("I102").Value = 5004.75
$ws.Range("K102").Value = 5004.75
$ws.Range("M102").Value = -3382.75
$ws.Range("H122").Value = 1311.1666
$ws.Range("I122").Value = 1248.5454
$ws.Range("K122").Value = 3745.6362
$ws.Range("M122").Value = -1295.6362
$ws.Range("H136").Value = 2433931.5
$ws.Range("I136").Value = 3257.4443
$ws.Range("K136").Value = 9772.332900000001
$ws.Range("M136").Value = -7222.332900000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 22151.574
$ws.Range("I20").Value = 9563.611000000001
$ws.Range("J20").Value = 32450.818
$ws.Range("K20").Value = 9563.611000000001
$ws.Range("L20").Value = 32450.818
$ws.Range("M20").Value = -9316.611000000001
$ws.Range("N20").Value = -32944.818
$ws.Range("H86").Value = 2934.4285
$ws.Range("I86").Value = 2994.261
$ws.Range("J86").Value = 2659.2
$ws.Range("K86").Value = 2994.261
$ws.Range("L86").Value = 2659.2
$ws.Range("M86").Value = -1871.261
$ws.Range("N86").Value = -4905.2
$ws.Range("H89").Value = 2934.4285
$ws.Range("I89").Value = 2994.261
$ws.Range("J89").Value = 2659.2
$ws.Range("K89").Value = 14971.305
$ws.Range("L89").Value = 13296
$ws.Range("M89").Value = -9355.305
$ws.Range("N89").Value = -24528
$ws.Range("H94").Value = 805.7273
$ws.Range("I94").Value = 557.8333
$ws.Range("J94").Value = 1103.2
$ws.Range("K94").Value = 557.8333
$ws.Range("L94").Value = 1103.2
$ws.Range("M94").Value = -106.8333
$ws.Range("N94").Value = -2005.2
$ws.Range("H99").Value = 1480.45
$ws.Range("J99").Value = 1541.5
$ws.Range("L99").Value = 1541.5
$ws.Range("N99").Value = -4537.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12568.206
$ws.Range("J58").Value = 100014
$ws.Range("L58").Value = 100014
$ws.Range("N58").Value = -100420
$ws.Range("H62").Value = 3856.7144
$ws.Range("I62").Value = 1999
$ws.Range("J62").Value = 5250
$ws.Range("K62").Value = 1999
$ws.Range("L62").Value = 5250
$ws.Range("M62").Value = -1375
$ws.Range("N62").Value = -6498
$ws.Range("H65").Value = 3856.7144
$ws.Range("I65").Value = 1999
$ws.Range("J65").Value = 5250
$ws.Range("K65").Value = 9995
$ws.Range("L65").Value = 26250
$ws.Range("M65").Value = -6875
$ws.Range("N65").Value = -32490
$ws.Range("H133").Value = 79999.39999999999
$ws.Range("J133").Value = 79999.39999999999
$ws.Range("L133").Value = 79999.39999999999
$ws.Range("N133").Value = -85059.39999999999
$ws.Range("H134").Value = 43485052
$ws.Range("I134").Value = 2224.875
$ws.Range("J134").Value = 142874370
$ws.Range("K134").Value = 6674.625
$ws.Range("L134").Value = 428623110
$ws.Range("M134").Value = -4139.625
$ws.Range("N134").Value = -428628180
$ws.Range("H136").Value = 12568.206
$ws.Range("J136").Value = 100014
$ws.Range("L136").Value = 300042
$ws.Range("N136").Value = -305142
$ws.Range("H137").Value = 53599.8
$ws.Range("J137").Value = 53599.8
$ws.Range("L137").Value = 53599.8
$ws.Range("N137").Value = -63799.8
$ws.Range("H138").Value = 64665.5
$ws.Range("J138").Value = 64665.5
$ws.Range("L138").Value = 64665.5
$ws.Range("N138").Value = -74945.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1037.5
$ws.Range("J68").Value = 1230
$ws.Range("L68").Value = 3690
$ws.Range("N68").Value = -5312
$ws.Range("H71").Value = 1037.5
$ws.Range("J71").Value = 1230
$ws.Range("L71").Value = 11070
$ws.Range("N71").Value = -19182
$ws.Range("H122").Value = 17934696
$ws.Range("J122").Value = 4722683.5
$ws.Range("L122").Value = 42504151.5
$ws.Range("N122").Value = -42509051.5
$ws.Range("H127").Value = 82426.55499999999
$ws.Range("J127").Value = 82426.55499999999
$ws.Range("L127").Value = 247279.665
$ws.Range("N127").Value = -257199.665

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 14929.363
$ws.Range("I70").Value = 14666.667
$ws.Range("K70").Value = 14666.667
$ws.Range("M70").Value = -14396.667
$ws.Range("H73").Value = 14929.363
$ws.Range("I73").Value = 14666.667
$ws.Range("K73").Value = 14666.667
$ws.Range("M73").Value = -13730.667
$ws.Range("H102").Value = 6663.3335
$ws.Range("I102").Value = 8062.6665
$ws.Range("J102").Value = 3165
$ws.Range("K102").Value = 8062.6665
$ws.Range("L102").Value = 3165
$ws.Range("M102").Value = -6440.6665
$ws.Range("N102").Value = -6409
$ws.Range("H122").Value = 2956.8667
$ws.Range("I122").Value = 2953.7856
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 8861.356800000001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -6411.356800000001
$ws.Range("N122").Value = -13900
$ws.Range("H126").Value = 5673.4
$ws.Range("I126").Value = 6823.727
$ws.Range("J126").Value = 2510
$ws.Range("K126").Value = 20471.181
$ws.Range("L126").Value = 7530
$ws.Range("M126").Value = -18001.181
$ws.Range("N126").Value = -12470

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4989.6665
$ws.Range("I40").Value = 4978
$ws.Range("K40").Value = 4978
$ws.Range("M40").Value = -4842
$ws.Range("H93").Value = 3948.353
$ws.Range("I93").Value = 4925.231
$ws.Range("J93").Value = 773.5
$ws.Range("K93").Value = 4925.231
$ws.Range("L93").Value = 773.5
$ws.Range("M93").Value = -3677.231
$ws.Range("N93").Value = -3269.5
$ws.Range("H100").Value = 2997.5
$ws.Range("I100").Value = 2667.8572
$ws.Range("K100").Value = 2667.8572
$ws.Range("M100").Value = -2126.8572
$ws.Range("H136").Value = 3551080.5
$ws.Range("I136").Value = 41811.832
$ws.Range("K136").Value = 125435.496
$ws.Range("M136").Value = -122885.496

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 616.3333
$ws.Range("I100").Value = 366.33334
$ws.Range("J100").Value = 866.3333
$ws.Range("K100").Value = 732.66668
$ws.Range("L100").Value = 1732.6666
$ws.Range("M100").Value = -191.66668
$ws.Range("N100").Value = -2814.6666
$ws.Range("H109").Value = 23000
$ws.Range("J109").Value = 23000
$ws.Range("L109").Value = 23000
$ws.Range("N109").Value = -25774
$ws.Range("H136").Value = 347934.88
$ws.Range("I136").Value = 1301.125
$ws.Range("J136").Value = 1734469.9
$ws.Range("K136").Value = 3903.375
$ws.Range("L136").Value = 5203409.699999999
$ws.Range("M136").Value = -1353.375
$ws.Range("N136").Value = -5208509.699999999
